$d = $word.ActiveDocument

# Locate the sentence that needs to be split/re-cased. Search without the
# smart quote (apostrophe) since Find matches both straight and curly quotes.
$rng = $d.Content
$found = $rng.Find.Execute("All Queries are based on Earnings before interests and taxes (EBIT) as an indicator of the company's profitability.", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not find target sentence"
}

# $rng.Start points right before "All"; the run actually begins one
# character earlier, at the leading space.
$baseStart = $rng.Start - 1

# Offsets (relative to $baseStart) bounding each of the 5 single-letter
# runs that must become their own <w:r> elements, plus their replacement
# (possibly re-cased) text.
$letterEdits = @(
    @{ From = 5;  To = 6;  Text = "q" },
    @{ From = 26; To = 27; Text = "E" },
    @{ From = 35; To = 36; Text = "B" },
    @{ From = 42; To = 43; Text = "I" },
    @{ From = 56; To = 57; Text = "T" }
)

foreach ($edit in $letterEdits) {
    $charRng = $d.Range($baseStart + $edit.From, $baseStart + $edit.To)
    # Toggling a character-level format and reverting it forces the run to
    # split from its neighbours even though the final formatting is
    # identical, mirroring how Word splits runs on manual retyping.
    $charRng.Font.Bold = 1
    $charRng.Text = $edit.Text
    $charRng2 = $d.Range($baseStart + $edit.From, $baseStart + $edit.To)
    $charRng2.Font.Bold = 0
}
